$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1380.8462
$ws.Range("J17").Value = 1455.7391
$ws.Range("L17").Value = 4367.2173
$ws.Range("N17").Value = -4703.2173
$ws.Range("H80").Value = 1398.2
$ws.Range("I80").Value = 997
$ws.Range("K80").Value = 2991
$ws.Range("M80").Value = -1993
$ws.Range("H83").Value = 1398.2
$ws.Range("I83").Value = 997
$ws.Range("K83").Value = 8973
$ws.Range("M83").Value = -3981
$ws.Range("H137").Value = 40744.28
$ws.Range("I137").Value = 82099.05
$ws.Range("J137").Value = 2835.75
$ws.Range("K137").Value = 246297.15
$ws.Range("L137").Value = 8507.25
$ws.Range("M137").Value = -243747.15
$ws.Range("N137").Value = -13607.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6257383
$ws.Range("I45").Value = 8992419
$ws.Range("J45").Value = 5871.4287
$ws.Range("K45").Value = 8992419
$ws.Range("L45").Value = 5871.4287
$ws.Range("M45").Value = -8992042
$ws.Range("N45").Value = -6625.4287
$ws.Range("H61").Value = 2726.72
$ws.Range("I61").Value = 2519.2727
$ws.Range("K61").Value = 2519.2727
$ws.Range("M61").Value = -2307.2727
$ws.Range("H74").Value = 42593.797
$ws.Range("I74").Value = 5364.7734
$ws.Range("J74").Value = 165914.94
$ws.Range("K74").Value = 5364.7734
$ws.Range("L74").Value = 165914.94
$ws.Range("M74").Value = -4490.7734
$ws.Range("N74").Value = -167662.94
$ws.Range("H77").Value = 42593.797
$ws.Range("I77").Value = 5364.7734
$ws.Range("J77").Value = 165914.94
$ws.Range("K77").Value = 26823.867
$ws.Range("L77").Value = 829574.7
$ws.Range("M77").Value = -22455.867
$ws.Range("N77").Value = -838310.7
$ws.Range("H132").Value = 2581.5
$ws.Range("I132").Value = 1643.4
$ws.Range("K132").Value = 4930.200000000001
$ws.Range("M132").Value = -2400.200000000001
$ws.Range("H136").Value = 2726.72
$ws.Range("I136").Value = 2519.2727
$ws.Range("K136").Value = 7557.8181
$ws.Range("M136").Value = -5007.8181
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2092.1667
$ws.Range("J64").Value = 2140
$ws.Range("L64").Value = 2140
$ws.Range("N64").Value = -2590
$ws.Range("H67").Value = 2092.1667
$ws.Range("J67").Value = 2140
$ws.Range("L67").Value = 2140
$ws.Range("N67").Value = -3700
$ws.Range("H86").Value = 5270204.5
$ws.Range("I86").Value = 10011339
$ws.Range("J86").Value = 2277.3333
$ws.Range("K86").Value = 10011339
$ws.Range("L86").Value = 2277.3333
$ws.Range("M86").Value = -10010216
$ws.Range("N86").Value = -4523.3333
$ws.Range("H89").Value = 5270204.5
$ws.Range("I89").Value = 10011339
$ws.Range("J89").Value = 2277.3333
$ws.Range("K89").Value = 50056695
$ws.Range("L89").Value = 11386.6665
$ws.Range("M89").Value = -50051079
$ws.Range("N89").Value = -22618.6665
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 5385.7144
$ws.Range("I32").Value = 2073.3333
$ws.Range("J32").Value = 7870
$ws.Range("K32").Value = 2073.3333
$ws.Range("L32").Value = 7870
$ws.Range("M32").Value = -1757.3333
$ws.Range("N32").Value = -8502
$ws.Range("H94").Value = 1036.1666
$ws.Range("I94").Value = 687.1667
$ws.Range("J94").Value = 1210.6666
$ws.Range("K94").Value = 687.1667
$ws.Range("L94").Value = 1210.6666
$ws.Range("M94").Value = -236.1667
$ws.Range("N94").Value = -2112.6666
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9621042
$ws.Range("I56").Value = 9621042
$ws.Range("K56").Value = 9621042
$ws.Range("M56").Value = -9620512
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H75").Value = 4851.3335
$ws.Range("J75").Value = 4851.3335
$ws.Range("L75").Value = 14554.0005
$ws.Range("N75").Value = -16550.0005
$ws.Range("H76").Value = 3416
$ws.Range("I76").Value = 3416
$ws.Range("K76").Value = 10248
$ws.Range("M76").Value = -9865
$ws.Range("H78").Value = 4851.3335
$ws.Range("J78").Value = 4851.3335
$ws.Range("L78").Value = 43662.0015
$ws.Range("N78").Value = -53646.0015
$ws.Range("H79").Value = 3416
$ws.Range("I79").Value = 3416
$ws.Range("K79").Value = 10248
$ws.Range("M79").Value = -8922
$ws.Range("H87").Value = 15600
$ws.Range("I87").Value = 11333.333
$ws.Range("J87").Value = 22000
$ws.Range("K87").Value = 33999.999
$ws.Range("L87").Value = 66000
$ws.Range("M87").Value = -32751.999
$ws.Range("N87").Value = -68496
$ws.Range("H90").Value = 15600
$ws.Range("I90").Value = 11333.333
$ws.Range("J90").Value = 22000
$ws.Range("K90").Value = 101999.997
$ws.Range("L90").Value = 198000
$ws.Range("M90").Value = -95759.997
$ws.Range("N90").Value = -210480
$ws.Range("H107").Value = 184.59259
$ws.Range("I107").Value = 174.38095
$ws.Range("J107").Value = 220.33333
$ws.Range("K107").Value = 523.1428500000001
$ws.Range("L107").Value = 660.99999
$ws.Range("M107").Value = 1396.85715
$ws.Range("N107").Value = -4500.99999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 945.4545000000001
$ws.Range("I2").Value = 1289.375
$ws.Range("J2").Value = 28.333334
$ws.Range("K2").Value = 1289.375
$ws.Range("L2").Value = 28.333334
$ws.Range("M2").Value = -1176.375
$ws.Range("N2").Value = -254.333334
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 4903.8887
$ws.Range("I36").Value = 4779.25
$ws.Range("J36").Value = 5003.6
$ws.Range("K36").Value = 4779.25
$ws.Range("L36").Value = 5003.6
$ws.Range("M36").Value = -4294.25
$ws.Range("N36").Value = -5973.6
$ws.Range("H43").Value = 9309.666999999999
$ws.Range("I43").Value = 1291.8334
$ws.Range("J43").Value = 25345.334
$ws.Range("K43").Value = 1291.8334
$ws.Range("L43").Value = 25345.334
$ws.Range("M43").Value = -1140.8334
$ws.Range("N43").Value = -25647.334
$ws.Range("H122").Value = 427626.53
$ws.Range("I122").Value = 686471.9
$ws.Range("J122").Value = 7002.875
$ws.Range("K122").Value = 2059415.7
$ws.Range("L122").Value = 21008.625
$ws.Range("M122").Value = -2056965.7
$ws.Range("N122").Value = -25908.625
$ws.Range("H132").Value = 3941.1765
$ws.Range("I132").Value = 3285.8572
$ws.Range("K132").Value = 9857.571599999999
$ws.Range("M132").Value = -7327.571599999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 910.3333
$ws.Range("I16").Value = 693.8570999999999
$ws.Range("J16").Value = 1668
$ws.Range("K16").Value = 693.8570999999999
$ws.Range("L16").Value = 1668
$ws.Range("M16").Value = -523.8570999999999
$ws.Range("N16").Value = -2008
$ws.Range("H41").Value = 32500
$ws.Range("I41").Value = 20000
$ws.Range("K41").Value = 20000
$ws.Range("M41").Value = -19562
$ws.Range("H122").Value = 4907.125
$ws.Range("I122").Value = 3027.1667
$ws.Range("K122").Value = 9081.500100000001
$ws.Range("M122").Value = -6631.500100000001
$ws.Range("H132").Value = 4298.84
$ws.Range("I132").Value = 3743.8333
$ws.Range("K132").Value = 11231.4999
$ws.Range("M132").Value = -8701.499899999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 14000
$ws.Range("J31").Value = 14000
$ws.Range("L31").Value = 14000
$ws.Range("N31").Value = -14696
$ws.Range("H113").Value = 892.7241
$ws.Range("I113").Value = 741.5625
$ws.Range("K113").Value = 2224.6875
$ws.Range("M113").Value = -54.6875
$ws.Range("H122").Value = 2112
$ws.Range("I122").Value = 1435.2727
$ws.Range("K122").Value = 4305.8181
$ws.Range("M122").Value = -1855.8181
$ws.Range("H132").Value = 37072320
$ws.Range("I132").Value = 40002068
$ws.Range("J132").Value = 450443.5
$ws.Range("K132").Value = 120006204
$ws.Range("L132").Value = 1351330.5
$ws.Range("M132").Value = -120003674
$ws.Range("N132").Value = -1356390.5
